# Update "想去人数" (interest counts) figures on the "展览" and "全部类型"
# sheets to reflect the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value  = 7736
$wsExhibition.Range("F5").Value  = 90
$wsExhibition.Range("F10").Value = 444
$wsExhibition.Range("F11").Value = 162
$wsExhibition.Range("F13").Value = 437
$wsExhibition.Range("F15").Value = 63
$wsExhibition.Range("F17").Value = 5592
$wsExhibition.Range("F18").Value = 153
$wsExhibition.Range("F19").Value = 217
$wsExhibition.Range("F20").Value = 1006
$wsExhibition.Range("F22").Value = 325

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 7737
$wsAll.Range("F5").Value  = 90
$wsAll.Range("F10").Value = 444
$wsAll.Range("F11").Value = 162
$wsAll.Range("F13").Value = 437
$wsAll.Range("F15").Value = 63
$wsAll.Range("F18").Value = 5592
$wsAll.Range("F20").Value = 153
$wsAll.Range("F21").Value = 217
$wsAll.Range("F22").Value = 1006
$wsAll.Range("F24").Value = 325
